$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.783.51'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.464.14'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.84%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.69'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.99%  '

$ws.Range("E9").Value = '  -0.69%  '

$ws.Range("E10").Value = '  -0.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.18'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.347'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.86'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000175'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.920.62'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.646.69'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.471.22'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.65'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -7.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.74'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.31'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.92%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.81'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.47%  '

$ws.Range("B22").Value = 'Polkadot'
$ws.Range("C22").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.09'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.64%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.84'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '646.13'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.65%  '

$ws.Range("E27").Value = '  -0.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0964'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.81%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.42'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.88'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.82'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.132'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.40%  '

$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.66'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.37'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.12%  '

$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '150.15'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.41%  '

$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.364'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.53'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.66'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.59%  '

$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.00%  '

$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₆0306'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.72'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.40'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.55'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.84%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.25'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.604'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0509'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0901'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.93%  '
